$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "61.276.44"
$ws.Range("E2").Value = "  +0.24%  "
$ws.Range("D3").Value = "2.923.43"
$ws.Range("E3").Value = "  -0.18%  "
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "596.62"
$ws.Range("E5").Value = "  +0.51%  "
$ws.Range("D6").Value = "145.05"
$ws.Range("E6").Value = "  -0.77%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("E8").Value = "  -1.09%  "
$ws.Range("D9").Value = "6.97"
$ws.Range("E9").Value = "  +0.77%  "
$ws.Range("E10").Value = "  -2.56%  "
$ws.Range("E11").Value = "  -0.80%  "
$ws.Range("E12").Value = "  -1.33%  "
$ws.Range("D13").Value = "33.39"
$ws.Range("E13").Value = "  -1.27%  "
$ws.Range("D15").Value = "3.408.81"
$ws.Range("E15").Value = "  -0.13%  "
$ws.Range("D16").Value = "61.309.26"
$ws.Range("E16").Value = "  +0.31%  "
$ws.Range("D17").Value = "2.923.59"
$ws.Range("E17").Value = "  -0.27%  "
$ws.Range("E18").Value = "  -0.68%  "
$ws.Range("D19").Value = "430.76"
$ws.Range("E19").Value = "  -0.28%  "
$ws.Range("E20").Value = "  -0.08%  "
$ws.Range("E21").Value = "  -1.23%  "
$ws.Range("E22").Value = "  -0.56%  "
$ws.Range("E23").Value = "  +0.53%  "
$ws.Range("D24").Value = "10.82"
$ws.Range("E24").Value = "  -2.18%  "
$ws.Range("D25").Value = "2.17"
$ws.Range("E25").Value = "  -2.44%  "
$ws.Range("D26").Value = "11.70"
$ws.Range("E26").Value = "  -2.87%  "
$ws.Range("E28").Value = "  -4.85%  "
$ws.Range("E29").Value = "  -0.61%  "
$ws.Range("E30").Value = "  -2.95%  "
$ws.Range("E31").Value = "  +1.19%  "
$ws.Range("D32").Value = "26.53"
$ws.Range("E32").Value = "  +0.18%  "
$ws.Range("E33").Value = "  -0.05%  "
$ws.Range("E34").Value = "  +2.70%  "
$ws.Range("E35").Value = "  -0.11%  "
$ws.Range("D36").Value = "5.60"
$ws.Range("E36").Value = "  -0.63%  "
$ws.Range("E37").Value = "  -3.36%  "
$ws.Range("E38").Value = "  -0.53%  "
$ws.Range("E39").Value = "  -1.93%  "
$ws.Range("E40").Value = "  -0.70%  "
$ws.Range("D41").Value = "42.15"
$ws.Range("E41").Value = "  +6.11%  "
$ws.Range("E42").Value = "  -2.44%  "
$ws.Range("E43").Value = "  -0.50%  "
$ws.Range("D44").Value = "2.697.28"
$ws.Range("E44").Value = "  -0.83%  "
$ws.Range("D45").Value = "133.76"
$ws.Range("E45").Value = "  +2.34%  "
$ws.Range("D46").Value = "360.48"
$ws.Range("E46").Value = "  -4.33%  "
$ws.Range("E47").Value = "  +0.04%  "
$ws.Range("D48").Value = "23.49"
$ws.Range("E48").Value = "  -2.86%  "
$ws.Range("E49").Value = "  -1.38%  "
$ws.Range("E50").Value = "  -2.11%  "
$ws.Range("E51").Value = "  -2.49%  "
